# Update the Wnt1-Fzd8 LR-pairs sheet with newly computed TPM-based values.
#
# The old sheet had six data rows: every combination of Sending cluster
# (ECs, FAPs) x Target cluster (ECs, FAPs, MuSCs). The new TPM numbers only
# keep the FAPs-as-sender combinations (what used to be rows 5-7), and the
# ECs-as-sender rows are dropped entirely. The surviving rows also get
# refreshed receptor / edge expression values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete "ECs" sending-cluster rows (old rows 2:4). This shifts
# the old "FAPs" sending-cluster rows (5:7) up to become rows 2:4.
$ws.Rows("2:4").Delete()

# Row 2 (Target cluster = ECs): refresh receptor-expressing-cell count /
# rate and every downstream edge-expression statistic.
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 3.235341333333333
$ws.Range("N2").Value = 9.706023999999999
$ws.Range("O2").Value = 0.2153734454473681
$ws.Range("P2").Value = 0.2153734454473681
$ws.Range("Q2").Value = 1.142744127875555
$ws.Range("R2").Value = 10.28469715088
$ws.Range("S2").Value = 0.2153734454473681
$ws.Range("T2").Value = 0.2153734454473681

# Row 3 (Target cluster = FAPs): receptor specificity and edge specificity
# values change; receptor counts / average & total expression stay as-is.
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.4841904166376352
$ws.Range("P3").Value = 0.4841904166376352
$ws.Range("S3").Value = 0.4841904166376352
$ws.Range("T3").Value = 0.4841904166376352

# Row 4 (Target cluster = MuSCs): same shape of update as row 3.
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("O4").Value = 0.3004361379149967
$ws.Range("P4").Value = 0.3004361379149967
$ws.Range("S4").Value = 0.3004361379149967
$ws.Range("T4").Value = 0.3004361379149967
